$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13954
$ws1.Range("F14").Value = 543
$ws1.Range("F18").Value = 14023
$ws1.Range("F19").Value = 371
$ws1.Range("F23").Value = 8306
$ws1.Range("F35").Value = 23
$ws1.Range("F45").Value = 5114

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13954
$ws4.Range("F14").Value = 543
$ws4.Range("F18").Value = 14023
$ws4.Range("F19").Value = 371
$ws4.Range("F23").Value = 8306
$ws4.Range("F35").Value = 23
$ws4.Range("F47").Value = 5114
